# Update "edit user" popup labels/values: Tenant/Sub-Tenant -> Partner/Sub Partner,
# UserName -> User Name, and refresh Role + Partner + Sub Partner columns with the
# new partner-based values (replacing the old tenant placeholder values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "User Name"
$ws.Range("D1").Value = "Partner"
$ws.Range("E1").Value = "Sub Partner"

# Row 2 - Jareen
$ws.Range("B2").Value = "Agent"
$ws.Range("D2").Value = "Altaworx"
$ws.Range("E2").Value = "Atlantech-AWX"

# Row 3 - Sindhuja
$ws.Range("B3").Value = "Agent Partner Admin"
$ws.Range("D3").Value = "Altaworx"
$ws.Range("E3").Value = "Castle Point-AWX"

# Row 4 - Anjani
$ws.Range("B4").Value = "Super Admin"
$ws.Range("D4").Value = "Altaworx"
$ws.Range("E4").Value = "CSV-AWX"

# Row 5 - Manoj
$ws.Range("B5").Value = "Agent"
$ws.Range("D5").Value = "Altaworx"
$ws.Range("E5").Value = "Frontier-AWX"

# Row 6 - Lohitha
$ws.Range("B6").Value = "Partner Admin"
$ws.Range("D6").Value = "Altaworx"
$ws.Range("E6").Value = "GoTech-AWX"

# Row 7 - Gopi
$ws.Range("B7").Value = "User"
$ws.Range("D7").Value = "Altaworx"
$ws.Range("E7").Value = "Local IT-AWX"

# Row 8 - Phani
$ws.Range("B8").Value = "Agent"
$ws.Range("D8").Value = "Altaworx"
$ws.Range("E8").Value = "Titanium-AWX"

# Row 9 - Pooja
$ws.Range("B9").Value = "Super Admin"
$ws.Range("D9").Value = "Altaworx"
$ws.Range("E9").Value = "Castle Point-AWX"

# Row 10 - Nikhil
$ws.Range("B10").Value = "Agent"
$ws.Range("D10").Value = "Altaworx"
$ws.Range("E10").Value = "Atlantech-AWX"

# Row 11 - Tejaswini
$ws.Range("B11").Value = "Agent Partner Admin"
$ws.Range("D11").Value = "Altaworx"
$ws.Range("E11").Value = "Frontier-AWX"

# Column widths widened to fit the new Partner / Sub Partner columns
$ws.Columns.Item(2).ColumnWidth = 17
$ws.Columns.Item(5).ColumnWidth = 15

# Restore the active selection left after editing
$ws.Range("E22").Select() | Out-Null
